$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037796343593383
$bf[0,2] = 1.044486559985052
$bf[0,3] = 1.046456344966288
$bf[0,4] = 1.057143380765936
$ws.Range("B2:F2").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037152754562994
$inn[0,1] = 1.042897267295881
$inn[0,2] = 1.047257473777299
$inn[0,3] = 1.049221729870916
$inn[0,4] = 1.059879157906199
$inn[0,5] = 1.04437830103743
$ws.Range("I2:N2").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.038722957236643
$bf[0,2] = 1.045174204342635
$bf[0,3] = 1.047255848676656
$bf[0,4] = 1.057969102210669
$ws.Range("B3:F3").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037295165172657
$inn[0,1] = 1.043468523255932
$inn[0,2] = 1.047756688465756
$inn[0,3] = 1.049832911356265
$inn[0,4] = 1.060518627133846
$inn[0,5] = 1.044950368246468
$ws.Range("I3:N3").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.039323247188092
$bf[0,2] = 1.045619594724897
$bf[0,3] = 1.047774126278463
$bf[0,4] = 1.058504248678144
$ws.Range("B4:F4").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037386175359054
$inn[0,1] = 1.043838242508672
$inn[0,2] = 1.048079473414449
$inn[0,3] = 1.05022867559141
$inn[0,4] = 1.060932603640604
$inn[0,5] = 1.045320612542954
$ws.Range("I4:N4").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.039575777467496
$bf[0,2] = 1.045806940139133
$bf[0,3] = 1.047992234786214
$bf[0,4] = 1.058729425816858
$ws.Range("B5:F5").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037424162729044
$inn[0,1] = 1.043993690045021
$inn[0,2] = 1.048215113411468
$inn[0,3] = 1.05039512270698
$inn[0,4] = 1.061106685075199
$inn[0,5] = 1.045476280832635
$ws.Range("I5:N5").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.039618188247749
$bf[0,2] = 1.045838402267072
$bf[0,3] = 1.048028869286182
$bf[0,4] = 1.058767245816647
$ws.Range("B6:F6").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037430524930064
$inn[0,1] = 1.044019791364857
$inn[0,2] = 1.048237884479439
$inn[0,3] = 1.050423073848697
$inn[0,4] = 1.061135916751314
$inn[0,5] = 1.045502419219339
$ws.Range("I6:N6").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.039326620850085
$bf[0,2] = 1.045622097639716
$bf[0,3] = 1.047777039774599
$bf[0,4] = 1.05850725671635
$ws.Range("B7:F7").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037386684022871
$inn[0,1] = 1.043840319538778
$inn[0,2] = 1.048081286075332
$inn[0,3] = 1.05023089940107
$inn[0,4] = 1.060934929547508
$inn[0,5] = 1.045322692522682
$ws.Range("I7:N7").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.038109349407852
$bf[0,2] = 1.044718860810558
$bf[0,3] = 1.046726344255241
$bf[0,4] = 1.057422260220113
$ws.Range("B8:F8").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.037201118351832
$inn[0,1] = 1.04309030894647
$inn[0,2] = 1.04742623472514
$inn[0,3] = 1.049428220879776
$inn[0,4] = 1.06009522780818
$inn[0,5] = 1.044571616829304
$ws.Range("I8:N8").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.035969850894696
$bf[0,2] = 1.043130674763824
$bf[0,3] = 1.044882198548658
$bf[0,4] = 1.05551694346111
$ws.Range("B9:F9").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036865432572956
$inn[0,1] = 1.041769349124878
$inn[0,2] = 1.046270165486244
$inn[0,3] = 1.04801607111638
$inn[0,4] = 1.058617141837636
$inn[0,5] = 1.04324878109325
$ws.Range("I9:N9").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034547288784873
$bf[0,2] = 1.042074291999803
$bf[0,3] = 1.043657786103207
$bf[0,4] = 1.054251267463764
$ws.Range("B10:F10").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036635833504705
$inn[0,1] = 1.04088922175804
$inn[0,2] = 1.045498327310359
$inn[0,3] = 1.047076245590065
$inn[0,4] = 1.057632899206667
$inn[0,5] = 1.042367403844599
$ws.Range("I10:N10").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033932214491339
$bf[0,2] = 1.041617458798077
$bf[0,3] = 1.043128813357303
$bf[0,4] = 1.053704314066608
$ws.Range("B11:F11").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036535044777636
$inn[0,1] = 1.040508253627085
$inn[0,2] = 1.04516386166381
$inn[0,3] = 1.046669688853882
$inn[0,4] = 1.057207002198779
$inn[0,5] = 1.0419858946952
$ws.Range("I11:N11").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033703885708178
$bf[0,2] = 1.041447860569683
$bf[0,3] = 1.042932512209539
$bf[0,4] = 1.053501317204002
$ws.Range("B12:F12").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.03649740202532
$inn[0,1] = 1.040366766265931
$inn[0,2] = 1.045039589230721
$inn[0,3] = 1.046518736106394
$inn[0,4] = 1.057048849630404
$inn[0,5] = 1.041844206405765
$ws.Range("I12:N12").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033752856777821
$bf[0,2] = 1.041484235867131
$bf[0,3] = 1.042974611168881
$bf[0,4] = 1.053544853179825
$ws.Range("B13:F13").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036505485805216
$inn[0,1] = 1.040397114807984
$inn[0,2] = 1.045066247743652
$inn[0,3] = 1.04655111322262
$inn[0,4] = 1.057082771853525
$inn[0,5] = 1.041874598046228
$ws.Range("I13:N13").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033913337956568
$bf[0,2] = 1.041603437910579
$bf[0,3] = 1.043112583302576
$bf[0,4] = 1.053687530886095
$ws.Range("B14:F14").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036531937401106
$inn[0,1] = 1.040496557798909
$inn[0,2] = 1.045153590013254
$inn[0,3] = 1.046657209810189
$inn[0,4] = 1.057193928340833
$inn[0,5] = 1.041974182257607
$ws.Range("I14:N14").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034012233926274
$bf[0,2] = 1.041676894195446
$bf[0,3] = 1.043197616772546
$bf[0,4] = 1.053775461378633
$ws.Range("B15:F15").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036548207912448
$inn[0,1] = 1.040557830748455
$inn[0,2] = 1.045207399602896
$inn[0,3] = 1.046722587479958
$inn[0,4] = 1.057262421468387
$inn[0,5] = 1.042035542221769
$ws.Range("I15:N15").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034588128407957
$bf[0,2] = 1.042104623037832
$bf[0,3] = 1.04369291783529
$bf[0,4] = 1.054287590136059
$ws.Range("B16:F16").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036642493691797
$inn[0,1] = 1.040914508259934
$inn[0,2] = 1.045520519436233
$inn[0,3] = 1.047103235835849
$inn[0,4] = 1.057661170753069
$inn[0,5] = 1.042392726256226
$ws.Range("I16:N16").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034949614879132
$bf[0,2] = 1.042373084530816
$bf[0,3] = 1.04400393117076
$bf[0,4] = 1.054609128798389
$ws.Range("B17:F17").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036701270099829
$inn[0,1] = 1.041138279054681
$inn[0,2] = 1.045716863950786
$inn[0,3] = 1.047342112905126
$inn[0,4] = 1.05791137341513
$inn[0,5] = 1.042616814831169
$ws.Range("I17:N17").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.035160551003916
$bf[0,2] = 1.04252973014976
$bf[0,3] = 1.044185456227075
$bf[0,4] = 1.054796782163976
$ws.Range("B18:F18").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036735421078879
$inn[0,1] = 1.041268813508045
$inn[0,2] = 1.045831363732669
$inn[0,3] = 1.047481483775407
$inn[0,4] = 1.058057339942244
$inn[0,5] = 1.042747534658429
$ws.Range("I18:N18").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.035232489519299
$bf[0,2] = 1.042583151801692
$bf[0,3] = 1.044247371244599
$bf[0,4] = 1.054860784918165
$ws.Range("B19:F19").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036747043220133
$inn[0,1] = 1.041313324506572
$inn[0,2] = 1.045870400979862
$inn[0,3] = 1.047529012031035
$inn[0,4] = 1.058107115382978
$inn[0,5] = 1.042792108867679
$ws.Range("I19:N19").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034910821774756
$bf[0,2] = 1.042344275289541
$bf[0,3] = 1.04397055034773
$bf[0,4] = 1.054574619837704
$ws.Range("B20:F20").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036694977628003
$inn[0,1] = 1.041114269240203
$inn[0,2] = 1.045695800583328
$inn[0,3] = 1.047316479724632
$inn[0,4] = 1.057884526177243
$inn[0,5] = 1.04259277092
$ws.Range("I20:N20").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033866076451525
$bf[0,2] = 1.04156833338277
$bf[0,3] = 1.043071948885913
$bf[0,4] = 1.053645511261196
$ws.Range("B21:F21").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036524153728589
$inn[0,1] = 1.040467273722649
$inn[0,2] = 1.045127870903109
$inn[0,3] = 1.046625965327919
$inn[0,4] = 1.057161194284217
$inn[0,5] = 1.0419448565946
$ws.Range("I21:N21").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033209998023382
$bf[0,2] = 1.041080989247351
$bf[0,3] = 1.042508021427965
$bf[0,4] = 1.053062305046366
$ws.Range("B22:F22").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036415562285513
$inn[0,1] = 1.040060605016918
$inn[0,2] = 1.044770577523307
$inn[0,3] = 1.046192162372326
$inn[0,4] = 1.056706665267738
$inn[0,5] = 1.041537610372661
$ws.Range("I22:N22").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033557721871787
$bf[0,2] = 1.041339289611274
$bf[0,3] = 1.042806869097761
$bf[0,4] = 1.053371381974298
$ws.Range("B23:F23").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036473241035455
$inn[0,1] = 1.040276175723877
$inn[0,2] = 1.044960005317728
$inn[0,3] = 1.046422095821186
$inn[0,4] = 1.056947594702358
$inn[0,5] = 1.04175348721475
$ws.Range("I23:N23").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034928350439543
$bf[0,2] = 1.042357292771
$bf[0,3] = 1.04398563334584
$bf[0,4] = 1.054590212626464
$ws.Range("B24:F24").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036697821334035
$inn[0,1] = 1.041125118202662
$inn[0,2] = 1.045705318288732
$inn[0,3] = 1.047328062138527
$inn[0,4] = 1.057896657203475
$inn[0,5] = 1.04260363528923
$ws.Range("I24:N24").Value = $inn

$bf = New-Object "object[,]" 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.03652230385677
$bf[0,2] = 1.043540841828093
$bf[0,3] = 1.045358077311575
$bf[0,4] = 1.056008722014586
$ws.Range("B25:F25").Value = $bf

$inn = New-Object "object[,]" 1,6
$inn[0,0] = 1.036953241956487
$inn[0,1] = 1.042110763868668
$inn[0,2] = 1.04656924092541
$inn[0,3] = 1.048380868167555
$inn[0,4] = 1.058999066628562
$inn[0,5] = 1.043590680685137
$ws.Range("I25:N25").Value = $inn

